$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 606.38
$ws.Range("C3").Value = 620.49
$ws.Range("C4").Value = 602.97
$ws.Range("C5").Value = 620.3
$ws.Range("C6").Value = 620.3
